$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the new label strings first, in the order they were originally
# authored, so the shared-strings table gets the same allocation order
# as the source workbook: Man Hours, Power, Up to 02/20, Original Budget,
# Pump, Tank, Budget Difference, %.
$ws.Cells.Item(15, 1).Value = "Man Hours"
$ws.Cells.Item(16, 1).Value = "Power"
$ws.Cells.Item(15, 5).Value = "Up to 02/20"
$ws.Cells.Item(18, 1).Value = "Original Budget"
$ws.Cells.Item(5, 1).Value = "Pump"
$ws.Cells.Item(2, 1).Value = "Tank"
$ws.Cells.Item(19, 1).Value = "Budget Difference"
$ws.Cells.Item(19, 6).Value = "%"

# Data rows: Label, Quantity, Price
$data = @(
    @(2, "Tank", 1, 750),
    @(3, "Heat Exchanger", 1, 600),
    @(4, "Flowmeter", 1, 60),
    @(5, "Pump", 1, 2000),
    @(6, "Thermocouple", 2, 10),
    @(7, "Electronic Control Valve", 1, 10),
    @(8, "Manual Valve", 2, 70),
    @(9, "Pipes 1 1/4", 5.5, 23),
    @(10, "Pipes 1", 6, 17),
    @(11, "Pipes 3/4", 2.5, 15),
    @(12, "Fittings 1 1/4", 5, 30),
    @(13, "Fittings 1", 10, 20),
    @(14, "Fittings 3/4", 6, 15)
)

foreach ($item in $data) {
    $row = $item[0]
    $ws.Cells.Item($row, 1).Value = $item[1]
    $ws.Cells.Item($row, 2).Value = $item[2]
    $ws.Cells.Item($row, 3).Value = $item[3]
}

# Rows 4 & 5 (Flowmeter, Pump) share a single formula, as in the source
# workbook (entered once across the D4:D5 range).
$ws.Range("D4:D5").Formula = "=C4*B4"

foreach ($item in $data) {
    $row = $item[0]
    if ($row -eq 4 -or $row -eq 5) { continue }
    $ws.Cells.Item($row, 4).Formula = "=C$row*B$row"
}

# Row 15: Man Hours
$ws.Cells.Item(15, 2).Formula = "=6+6+9+9+7+7+3+1"
$ws.Cells.Item(15, 3).Value = 40
$ws.Cells.Item(15, 4).Formula = "=C15*B15"

# Row 16: Power
$ws.Cells.Item(16, 3).Value = 1
$ws.Cells.Item(16, 4).Formula = "=C16*B16"

# Row 17: Total
$ws.Cells.Item(17, 1).Value = "Total"
$ws.Cells.Item(17, 4).Formula = "=SUM(D2:D16)"

# Row 18: Original Budget
$ws.Cells.Item(18, 4).Value = 5530

# Row 19: Budget Difference
$ws.Cells.Item(19, 4).Formula = "=D18-D17"
$ws.Cells.Item(19, 5).Formula = "=D17/D18*100"

$ws.Range("A20").Select()
